$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "497×4=1988" "693×6=4158"
Replace-Text "678×4=2712" "714×2=1428"
Replace-Text "741×6=4446" "832×4=3328"
Replace-Text "289×7=2023" "448×9=4032"
Replace-Text "729×6=4374" "313×7=2191"
Replace-Text "366×7=2562" "274×2=548"
Replace-Text "535×5=2675" "620×8=4960"
Replace-Text "929×3=2787" "781×2=1562"
Replace-Text "734×3=2202" "529×8=4232"
Replace-Text "716×4=2864" "955×9=8595"
Replace-Text "294×2=588" "461×6=2766"
Replace-Text "960×9=8640" "673×4=2692"
Replace-Text "169×9=1521" "302×6=1812"
Replace-Text "503×9=4527" "791×7=5537"
Replace-Text "673×2=1346" "843×7=5901"
Replace-Text "968×8=7744" "674×8=5392"
Replace-Text "365×6=2190" "211×9=1899"
Replace-Text "264×8=2112" "526×8=4208"
Replace-Text "617×5=3085" "305×2=610"
Replace-Text "824×5=4120" "439×3=1317"
Replace-Text "659×6=3954" "117×4=468"
Replace-Text "154×2=308" "771×2=1542"
Replace-Text "903×8=7224" "468×2=936"
Replace-Text "397×4=1588" "525×5=2625"
Replace-Text "809×5=4045" "431×7=3017"
